# Apply the "cibmtr-reporting-ig" update to the ValueSet workbook.
#
# Summary of changes (per the target diff):
#  - Metadata sheet: Version 0.1.6 -> 0.1.7
#  - Metadata sheet: Status active -> draft
#  - Metadata sheet: Date -> 2024-08-27T12:23:18-05:00
#  - Metadata sheet: first Contact row value updated to the publisher contact text
#  - Metadata sheet: second Contact row value updated to "Bob Milius (bmilius@nmdp.org)"
#  - Metadata sheet: a new "Jurisdiction" row (blank value) is inserted right
#    after the two Contact rows, pushing Description/Purpose/Copyright/
#    Immutable down by one row.
#  - Include from RxNorm sheet: content is unchanged (only shared-string
#    reindexing happens as a side effect of the Metadata sheet edits).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update simple scalar metadata values -------------------------------
$ws.Range("B3").Value = "0.1.7"
$ws.Range("B6").Value = "draft"
$ws.Range("B8").Value = "2024-08-27T12:23:18-05:00"

# Row 10 is the first "Contact" row -> update its value to the publisher contact.
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# Row 11 is the second "Contact" row -> update its value to the named contact.
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# --- Insert a new "Jurisdiction" row after the Contact rows (row 12) ----
# Extend the table's shared row style one row further (row 16 doesn't exist
# yet) by copying the format from the last existing data row, then shift the
# existing rows 12-15 down into 13-16 and write the new row 12 values. This
# keeps every cell on the same shared style used throughout the table
# instead of generating extra/unused styles the way Rows.Insert() does.
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($r = 15; $r -ge 12; $r--) {
  $dst = $r + 1
  $ws.Cells.Item($dst, 1).Value = $ws.Cells.Item($r, 1).Value2
  $ws.Cells.Item($dst, 2).Value = $ws.Cells.Item($r, 2).Value2
}

$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""
